$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws 'D2' '30.163.30'
Set-TextCell $ws 'E2' '  +0.68%  '
Set-TextCell $ws 'D3' '1.922.57'
Set-TextCell $ws 'E3' '  +2.94%  '
Set-TextCell $ws 'E4' '  +0.09%  '
Set-TextCell $ws 'D5' '319.74'
Set-TextCell $ws 'E5' '  +0.17%  '
Set-TextCell $ws 'E6' '  +0.08%  '
Set-TextCell $ws 'D7' '0.5076'
Set-TextCell $ws 'E7' '  -0.25%  '
Set-TextCell $ws 'D8' '0.4076'
Set-TextCell $ws 'E8' '  +3.51%  '
Set-TextCell $ws 'D9' '0.08345'
Set-TextCell $ws 'E9' '  +1.77%  '
Set-TextCell $ws 'D10' '1.120'
Set-TextCell $ws 'E10' '  +2.45%  '
Set-TextCell $ws 'B11' 'OKB'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D11' '42.04'
Set-TextCell $ws 'E11' '  -0.43%  '
Set-TextCell $ws 'B12' 'Solana'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws 'D12' '24.05'
Set-TextCell $ws 'E12' '  +4.76%  '
Set-TextCell $ws 'B13' 'WrappedEther'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D13' '1.922.61'
Set-TextCell $ws 'E13' '  +3.42%  '
Set-TextCell $ws 'B14' 'Polkadot'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws 'D14' '6.432'
Set-TextCell $ws 'E14' '  +2.42%  '
Set-TextCell $ws 'B15' 'Chainlink'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D15' '7.256'
Set-TextCell $ws 'E15' '  +1.09%  '
Set-TextCell $ws 'B16' 'BinanceUSD'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell $ws 'D16' '1.004'
Set-TextCell $ws 'E16' '  +0.24%  '
Set-TextCell $ws 'B17' 'Litecoin'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws 'D17' '92.74'
Set-TextCell $ws 'E17' '  +0.79%  '
Set-TextCell $ws 'B18' 'ShibaInu'
Set-TextCell $ws 'C18' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D18' '0.00001096'
Set-TextCell $ws 'E18' '  +0.82%  '
Set-TextCell $ws 'B19' 'TRON'
Set-TextCell $ws 'C19' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D19' '0.06512'
Set-TextCell $ws 'E19' '  +1.92%  '
Set-TextCell $ws 'B20' 'Avalanche'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D20' '18.52'
Set-TextCell $ws 'E20' '  +3.55%  '
Set-TextCell $ws 'B21' 'Dai'
Set-TextCell $ws 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D21' '1.001'
Set-TextCell $ws 'E21' '  +0.07%  '
Set-TextCell $ws 'B22' 'Uniswap'
Set-TextCell $ws 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws 'D22' '5.961'
Set-TextCell $ws 'E22' '  +2.38%  '
Set-TextCell $ws 'B23' 'WrappedBTC'
Set-TextCell $ws 'C23' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D23' '30.170.30'
Set-TextCell $ws 'E23' '  +0.72%  '
Set-TextCell $ws 'B24' 'Cosmos'
Set-TextCell $ws 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D24' '11.37'
Set-TextCell $ws 'E24' '  +2.34%  '
Set-TextCell $ws 'B25' 'Toncoin'
Set-TextCell $ws 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D25' '2.195'
Set-TextCell $ws 'E25' '  +1.11%  '
Set-TextCell $ws 'B26' 'WrappedliquidstakedEther2.0'
Set-TextCell $ws 'C26' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws 'D26' '2.143.57'
Set-TextCell $ws 'E26' '  +3.08%  '
Set-TextCell $ws 'B27' 'EthereumClassic'
Set-TextCell $ws 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D27' '21.98'
Set-TextCell $ws 'E27' '  +4.87%  '
Set-TextCell $ws 'B28' 'Monero'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D28' '162.79'
Set-TextCell $ws 'E28' '  +0.94%  '
Set-TextCell $ws 'B29' 'LidoDAOToken'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D29' '2.267'
Set-TextCell $ws 'E29' '  +1.99%  '
Set-TextCell $ws 'B30' 'BitcoinCash'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws 'D30' '129.00'
Set-TextCell $ws 'E30' '  +1.16%  '
Set-TextCell $ws 'B31' 'ImmutableX'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D31' '1.138'
Set-TextCell $ws 'E31' '  +7.23%  '
Set-TextCell $ws 'B32' 'Stellar'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D32' '0.1047'
Set-TextCell $ws 'E32' '  +1.26%  '
Set-TextCell $ws 'B33' 'Filecoin'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D33' '5.962'
Set-TextCell $ws 'E33' '  +0.63%  '
Set-TextCell $ws 'B34' 'HuobiToken'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws 'D34' '3.800'
Set-TextCell $ws 'E34' '  +1.88%  '
Set-TextCell $ws 'B35' 'VeChain'
Set-TextCell $ws 'C35' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D35' '0.02453'
Set-TextCell $ws 'E35' '  +1.07%  '
Set-TextCell $ws 'B36' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D36' '5.321'
Set-TextCell $ws 'E36' '  +2.17%  '
Set-TextCell $ws 'B37' 'Hedera'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D37' '0.06445'
Set-TextCell $ws 'E37' '  +1.61%  '
Set-TextCell $ws 'B38' 'ARBITRUM'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D38' '1.217'
Set-TextCell $ws 'E38' '  +3.85%  '
Set-TextCell $ws 'B39' 'Algorand'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws 'D39' '0.2147'
Set-TextCell $ws 'E39' '  +0.39%  '
Set-TextCell $ws 'B40' 'TheSandbox'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws 'D40' '0.6508'
Set-TextCell $ws 'E40' '  +3.30%  '
Set-TextCell $ws 'B41' 'FraxShare'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws 'D41' '8.596'
Set-TextCell $ws 'E41' '  +1.20%  '
Set-TextCell $ws 'B42' 'Aptos'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D42' '11.45'
Set-TextCell $ws 'E42' '  +1.66%  '
Set-TextCell $ws 'B43' 'TrustWalletToken'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws 'D43' '1.212'
Set-TextCell $ws 'E43' '  +0.82%  '
Set-TextCell $ws 'B44' 'EnergySwap'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D44' '13.34'
Set-TextCell $ws 'E44' '  +3.25%  '
Set-TextCell $ws 'D45' '0.6058'
Set-TextCell $ws 'E45' '  +2.78%  '
Set-TextCell $ws 'B46' 'NEARProtocol'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D46' '2.184'
Set-TextCell $ws 'E46' '  +9.09%  '
Set-TextCell $ws 'B47' 'PancakeSwap'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D47' '3.626'
Set-TextCell $ws 'E47' '  -0.16%  '
Set-TextCell $ws 'B48' 'Quant'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws 'D48' '122.87'
Set-TextCell $ws 'E48' '  +0.37%  '
Set-TextCell $ws 'B49' 'EOS'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws 'D49' '1.212'
Set-TextCell $ws 'E49' '  +0.82%  '
Set-TextCell $ws 'B50' 'WEMIXTOKEN'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws 'D50' '1.132'
Set-TextCell $ws 'E50' '  +1.41%  '
Set-TextCell $ws 'B51' 'Aave'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D51' '78.14'
Set-TextCell $ws 'E51' '  +1.88%  '
